$d = $word.ActiveDocument

# Locate the target paragraph robustly (search by a distinctive text fragment
# rather than a hard-coded paragraph index).
$paras = $d.Paragraphs
$target = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $t = $paras.Item($i).Range.Text
    if ($t -like '*I wanted to warn you that I won''t spend the night you tomorrow*') {
        $target = $paras.Item($i)
        break
    }
}

if ($target -eq $null) {
    throw "Could not find target paragraph"
}

$pStart = $target.Range.Start
$full = $target.Range.Text

# Find the insertion point: right before "spend the night..."
$marker = "spend the night"
$idx = $full.IndexOf($marker)
if ($idx -lt 0) {
    throw "Could not find insertion marker in paragraph"
}
$insertPos = $pStart + $idx

# Insert the new text "be able to " right before "spend the night ..."
$insRange = $d.Range($insertPos, $insertPos)
$insRange.InsertBefore("be able to ")

$insertedLen = "be able to ".Length

# Force the paragraph's single (re-merged) run to split into three runs at
# our two boundaries by nudging (and restoring) direct character formatting
# over each sub-range in turn. This makes the engine materialize separate
# <w:r> runs at those boundaries while keeping the run formatting identical
# to the surrounding text.
$rngBefore = $d.Range($pStart, $insertPos)
$rngBefore.Bold = 1
$rngBefore.Bold = 0

$rngMiddle = $d.Range($insertPos, $insertPos + $insertedLen)
$rngMiddle.Bold = 1
$rngMiddle.Bold = 0

Write-Output "Final paragraph text: $($target.Range.Text)"
